$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------
# Sheet 1 (position 1, physical sheet1.xml): GNG_TO... -> RS_TO...
# Shrinks from 4 data rows (A1:B5) to 2 data rows (A1:B3)
# --------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "RS_TO-16515889675295417"
$ws1.Rows.Item(5).Delete()
$ws1.Rows.Item(4).Delete()
$ws1.Range("B2").Value = "eyes closed"
$ws1.Range("B3").Value = "eyes open"

# --------------------------------------------------------------------
# Sheet 2 (position 2, physical sheet2.xml): NB_TO... -> GNG_TO...
# Shrinks from 9 data rows (A1:B10) to 4 data rows (A1:B5)
# --------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "GNG_TO-16515889675695746"
$ws2.Rows.Item(10).Delete()
$ws2.Rows.Item(9).Delete()
$ws2.Rows.Item(8).Delete()
$ws2.Rows.Item(7).Delete()
$ws2.Rows.Item(6).Delete()
$ws2.Range("B2").Value = "go_stims-16515889675331328.csv"
$ws2.Range("B3").Value = "GNG_stims-16515889675514083.csv"
$ws2.Range("B4").Value = "go_stims-16515889675534108.csv"
$ws2.Range("B5").Value = "GNG_stims-16515889675684946.csv"

# --------------------------------------------------------------------
# Sheet 3 (position 3, physical sheet3.xml): RS_TO... -> vSAT_TO...
# Grows from 2 data rows (A1:B3) to 4 data rows (A1:B5)
# --------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "vSAT_TO-1651588967645024"
$ws3.Range("A2").Copy($ws3.Range("A4:A5"))
$ws3.Range("A4").Value = 2
$ws3.Range("A5").Value = 3
$ws3.Range("B2").Value = "vSAT_stims-16515889676306102.csv"
$ws3.Range("B3").Value = "vSAT_stims-1651588967613811.csv"
$ws3.Range("B4").Value = "SAT_stims-16515889675980754.csv"
$ws3.Range("B5").Value = "SAT_stims-1651588967576935.csv"

# --------------------------------------------------------------------
# Sheet 4 (position 4, physical sheet4.xml): TOL_TO... -> TOL_TO...
# Same row count (A1:B7), only refreshed filenames
# --------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16515889677121215"
$ws4.Range("B2").Value = "MM_stims-1651588967676288.csv"
$ws4.Range("B3").Value = "ZM_stims-16515889676522746.csv"
$ws4.Range("B4").Value = "MM_stims-16515889676922998.csv"
$ws4.Range("B5").Value = "ZM_stims-16515889676782758.csv"
$ws4.Range("B6").Value = "MM_stims-1651588967711121.csv"
$ws4.Range("B7").Value = "ZM_stims-16515889676932728.csv"

# --------------------------------------------------------------------
# Sheet 5 (position 5, physical sheet5.xml): vSAT_TO... -> NB_TO...
# Grows from 4 data rows (A1:B5) to 9 data rows (A1:B10)
# --------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "NB_TO-165158896979931"
$ws5.Range("A2").Copy($ws5.Range("A6:A10"))
$ws5.Range("A6").Value = 4
$ws5.Range("A7").Value = 5
$ws5.Range("A8").Value = 6
$ws5.Range("A9").Value = 7
$ws5.Range("A10").Value = 8
$ws5.Range("B2").Value = "TB-16515889697747633.csv"
$ws5.Range("B3").Value = "OB-16515889685906258.csv"
$ws5.Range("B4").Value = "OB-16515889689589708.csv"
$ws5.Range("B5").Value = "ZB-match_9-16515889679291706.csv"
$ws5.Range("B6").Value = "TB-16515889689775045.csv"
$ws5.Range("B7").Value = "OB-1651588968415669.csv"
$ws5.Range("B8").Value = "TB-16515889696432564.csv"
$ws5.Range("B9").Value = "ZB-match_0-16515889679909225.csv"
$ws5.Range("B10").Value = "ZB-match_5-16515889680279288.csv"
